$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (word change + capitalization of header cells)
$ws.Range("B8").Value = "küszöbérték: 0.86; LLM modell: gemini-2.0-flash; top_k=50"
$ws.Range("B9").Value = " Kérdések száma"
$ws.Range("C9").Value = "Embedding  generálásai idő átlaga"
$ws.Range("D9").Value = "Kontextus összeállitási idő átlaga"
$ws.Range("F9").Value = "Teljes feldoldozási idő átlaga"
$ws.Range("G9").Value = "Szemantikus hasonlóság mérékének  (BERTScore F1) átlaga (0-1) között"

$ws.Range("E20").Select()
